$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the three placeholder name labels (B1:B3)
$ws.Range("B1").Value = "direc"
$ws.Range("B2").Value = "jefe"
$ws.Range("B3").Value = "coo"

# Reset the "Periodo Enero-Julio" counts (B4:B12) to 0
$ws.Range("B4").Value = 0
$ws.Range("B5").Value = 0
$ws.Range("B6").Value = 0
$ws.Range("B7").Value = 0
$ws.Range("B8").Value = 0
$ws.Range("B9").Value = 0
$ws.Range("B10").Value = 0
$ws.Range("B11").Value = 0
$ws.Range("B12").Value = 0

# Update the "Periodo Agosto-Diciembre" counts (B13:B21)
$ws.Range("B13").Value = 80
$ws.Range("B14").Value = 10
$ws.Range("B15").Value = 10
$ws.Range("B16").Value = 10
$ws.Range("B17").Value = 10
$ws.Range("B18").Value = 10
$ws.Range("B19").Value = 10
$ws.Range("B20").Value = 10
$ws.Range("B21").Value = 10

$wb.Save()
